# Update NATMI LR-pair sheet (Tnf -> Tnfrsf21) with newly recomputed TPM
# values, and drop the rows where "MuSCs" is the sending cluster (the
# new TPM run no longer produces ligand/receptor-expressing edges from
# that sender).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-10 (Sending cluster = MuSCs) are removed entirely; deleting the
# range shifts nothing else up from below since they were the last rows.
$ws.Rows("8:10").Delete()

# Sending/ligand/receptor/target-cluster labels are unchanged text-wise
# (still ECs/FAPs/MuSCs clusters signalling through Tnf -> Tnfrsf21);
# re-assert them defensively so the text content is correct either way.
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Tnf"
$ws.Range("C2").Value2 = "Tnfrsf21"
$ws.Range("D2").Value2 = "ECs"

$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Tnf"
$ws.Range("C3").Value2 = "Tnfrsf21"
$ws.Range("D3").Value2 = "FAPs"

$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Tnf"
$ws.Range("C4").Value2 = "Tnfrsf21"
$ws.Range("D4").Value2 = "MuSCs"

$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Tnf"
$ws.Range("C5").Value2 = "Tnfrsf21"
$ws.Range("D5").Value2 = "ECs"

$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Tnf"
$ws.Range("C6").Value2 = "Tnfrsf21"
$ws.Range("D6").Value2 = "FAPs"

$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Tnf"
$ws.Range("C7").Value2 = "Tnfrsf21"
$ws.Range("D7").Value2 = "MuSCs"

# Recomputed TPM-derived numeric columns (G,H,I,J,M,N,O,P,Q,R,S,T).
# E, F, K, L are untouched by the new run.

# Row 2: ECs -> Tnf/Tnfrsf21 -> ECs
$ws.Range("G2").Value2 = 2.913576333333333
$ws.Range("H2").Value2 = 8.740729
$ws.Range("I2").Value2 = 0.8649322955011439
$ws.Range("J2").Value2 = 0.8649322955011439
$ws.Range("M2").Value2 = 17.46627766666667
$ws.Range("N2").Value2 = 52.398833
$ws.Range("O2").Value2 = 0.2609791297364465
$ws.Range("P2").Value2 = 0.2609791297364465
$ws.Range("Q2").Value2 = 50.88933324102855
$ws.Range("R2").Value2 = 458.003999169257
$ws.Range("S2").Value2 = 0.2257292777608355
$ws.Range("T2").Value2 = 0.2257292777608355

# Row 3: ECs -> Tnf/Tnfrsf21 -> FAPs
$ws.Range("G3").Value2 = 2.913576333333333
$ws.Range("H3").Value2 = 8.740729
$ws.Range("I3").Value2 = 0.8649322955011439
$ws.Range("J3").Value2 = 0.8649322955011439
$ws.Range("M3").Value2 = 2.488505666666667
$ws.Range("N3").Value2 = 7.465517
$ws.Range("O3").Value2 = 0.03718296798122674
$ws.Range("P3").Value2 = 0.03718296798122674
$ws.Range("Q3").Value2 = 7.250451215765889
$ws.Range("R3").Value2 = 65.254060941893
$ws.Range("S3").Value2 = 0.03216074984954798
$ws.Range("T3").Value2 = 0.03216074984954798

# Row 4: ECs -> Tnf/Tnfrsf21 -> MuSCs
$ws.Range("G4").Value2 = 2.913576333333333
$ws.Range("H4").Value2 = 8.740729
$ws.Range("I4").Value2 = 0.8649322955011439
$ws.Range("J4").Value2 = 0.8649322955011439
$ws.Range("M4").Value2 = 46.97117233333334
$ws.Range("N4").Value2 = 140.913517
$ws.Range("O4").Value2 = 0.7018379022823268
$ws.Range("P4").Value2 = 0.7018379022823268
$ws.Range("Q4").Value2 = 136.8540960593214
$ws.Range("R4").Value2 = 1231.686864533893
$ws.Range("S4").Value2 = 0.6070422678907604
$ws.Range("T4").Value2 = 0.6070422678907604

# Row 5: FAPs -> Tnf/Tnfrsf21 -> ECs
$ws.Range("G5").Value2 = 0.4549836666666667
$ws.Range("H5").Value2 = 1.364951
$ws.Range("I5").Value2 = 0.1350677044988561
$ws.Range("J5").Value2 = 0.1350677044988561
$ws.Range("M5").Value2 = 17.46627766666667
$ws.Range("N5").Value2 = 52.398833
$ws.Range("O5").Value2 = 0.2609791297364465
$ws.Range("P5").Value2 = 0.2609791297364465
$ws.Range("Q5").Value2 = 7.946871055798111
$ws.Range("R5").Value2 = 71.521839502183
$ws.Range("S5").Value2 = 0.03524985197561098
$ws.Range("T5").Value2 = 0.03524985197561099

# Row 6: FAPs -> Tnf/Tnfrsf21 -> FAPs
$ws.Range("G6").Value2 = 0.4549836666666667
$ws.Range("H6").Value2 = 1.364951
$ws.Range("I6").Value2 = 0.1350677044988561
$ws.Range("J6").Value2 = 0.1350677044988561
$ws.Range("M6").Value2 = 2.488505666666667
$ws.Range("N6").Value2 = 7.465517
$ws.Range("O6").Value2 = 0.03718296798122674
$ws.Range("P6").Value2 = 0.03718296798122674
$ws.Range("Q6").Value2 = 1.132229432740778
$ws.Range("R6").Value2 = 10.190064894667
$ws.Range("S6").Value2 = 0.00502221813167876
$ws.Range("T6").Value2 = 0.005022218131678761

# Row 7: FAPs -> Tnf/Tnfrsf21 -> MuSCs
$ws.Range("G7").Value2 = 0.4549836666666667
$ws.Range("H7").Value2 = 1.364951
$ws.Range("I7").Value2 = 0.1350677044988561
$ws.Range("J7").Value2 = 0.1350677044988561
$ws.Range("M7").Value2 = 46.97117233333334
$ws.Range("N7").Value2 = 140.913517
$ws.Range("O7").Value2 = 0.7018379022823268
$ws.Range("P7").Value2 = 0.7018379022823268
$ws.Range("Q7").Value2 = 21.37111621585189
$ws.Range("R7").Value2 = 192.340045942667
$ws.Range("S7").Value2 = 0.09479563439156634
$ws.Range("T7").Value2 = 0.09479563439156637

Write-Output "Tnf-Tnfrsf21 sheet updated: 6 data rows (ECs/FAPs senders), new TPM values applied."
